$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 "Desarrollo de integración con pasarela de pagos." -> progress 50% -> 100%
$ws.Range("E30").Value = 100

# Row 31 "Desarrollo de servicios web para los datos de los reportes." becomes the new
# sales-report activity, now fully complete.
$ws.Range("C31").Value = "Desarrollo de servicios web para el reporte de ventas"
$ws.Range("E31").Value = 100

# Row 32 ("Desarrollo de servicios web para la gestión de la contabilidad.") is removed
# entirely; everything below shifts up one row.
$ws.Rows(32).Delete()

# Update the view's current selection to match the author's final cursor position.
$ws.Activate()
$excel.Goto($ws.Range("F27"), $true)
